# Updates cryptos price (D) and volume/change (E) columns per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values are plain numeric-looking strings (e.g. '228.35'); without
# forcing a Text number format first, Excel's COM layer auto-converts them to
# numbers. Force those specific cells to Text so the value round-trips as a string,
# matching the source data (which stores all Price/Volume cells as text).
$textForceCells = @("D5", "D7", "D10", "D13", "D16", "D19", "D22", "D24", "D26", "D27", "D29", "D30", "D32", "D35", "D37", "D38", "D39", "D42", "D46", "D47", "D48", "D49")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "38.852.69"
$ws.Range("E2").Value = "  +2.86%  "
$ws.Range("D3").Value = "2.092.38"
$ws.Range("E3").Value = "  +2.38%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "228.35"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").Value = "60.59"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("D10").Value = "0.0836"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "2.403.34"
$ws.Range("E12").Value = "  +2.46%  "
$ws.Range("D13").Value = "14.97"
$ws.Range("E13").Value = "  +4.13%  "
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("E15").Value = "  +4.07%  "
$ws.Range("D16").Value = "5.46"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "2.106.49"
$ws.Range("E17").Value = "  +3.13%  "
$ws.Range("D18").Value = "38.716.96"
$ws.Range("E18").Value = "  +2.67%  "
$ws.Range("D19").Value = "71.61"
$ws.Range("E19").Value = "  +3.33%  "
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("D21").Value = "0.0₃0837"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").Value = "227.21"
$ws.Range("E22").Value = "  +2.24%  "
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").Value = "171.30"
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("D27").Value = "9.49"
$ws.Range("E27").Value = "  +1.96%  "
$ws.Range("E28").Value = "  +9.40%  "
$ws.Range("D29").Value = "1.49"
$ws.Range("E29").Value = "  +15.98%  "
$ws.Range("D30").Value = "19.17"
$ws.Range("E30").Value = "  +2.25%  "
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").Value = "2.39"
$ws.Range("E32").Value = "  +6.10%  "
$ws.Range("E33").Value = "  +2.94%  "
$ws.Range("E34").Value = "  +4.53%  "
$ws.Range("D35").Value = "0.0611"
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").Value = "2.38"
$ws.Range("E37").Value = "  +1.80%  "
$ws.Range("D38").Value = "3.57"
$ws.Range("E38").Value = "  +3.26%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").Value = "1.541.35"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("D42").Value = "100.93"
$ws.Range("E42").Value = "  +3.38%  "
$ws.Range("E43").Value = "  +4.40%  "
$ws.Range("E44").Value = "  -0.72%  "
$ws.Range("E45").Value = "  +3.95%  "
$ws.Range("D46").Value = "7.65"
$ws.Range("E46").Value = "  +9.16%  "
$ws.Range("D47").Value = "1.13"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("D48").Value = "4.12"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").Value = "1.04"
$ws.Range("E49").Value = "  +3.03%  "
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("D51").Value = "2.288.79"
$ws.Range("E51").Value = "  +2.44%  "
